# Weekly data refresh: insert two new daily-price records at the top of the
# "Cebollín" price table (rows 167-168), pushing the existing 121 records
# (old rows 167-287) down by two rows to 169-289.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 167 (each Insert() pushes rows down).
$ws.Rows.Item(167).Insert()
$ws.Rows.Item(167).Insert()

# --- New row 167 ---
$ws.Range("A167").Value = 10
$ws.Range("B167").Value = "Vega Modelo de Temuco"
$ws.Range("C167").Value = "La Araucanía"
$ws.Range("D167").Value = (Get-Date -Year 2022 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E167").Value = 9
$ws.Range("F167").Value = 100112037
$ws.Range("G167").Value = "Cebollín"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 185
$ws.Range("K167").Value = 8000
$ws.Range("L167").Value = 8000
$ws.Range("M167").Value = 8000
$ws.Range("N167").Value = "`$/docena de paquetes"
$ws.Range("O167").Value = "Provincia de Cautín"
$ws.Range("P167").Value = 667
$ws.Range("Q167").Value = 12
$ws.Range("R167").Value = "Hortaliza"

# --- New row 168 ---
$ws.Range("A168").Value = 10
$ws.Range("B168").Value = "Vega Modelo de Temuco"
$ws.Range("C168").Value = "La Araucanía"
$ws.Range("D168").Value = (Get-Date -Year 2022 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112037
$ws.Range("G168").Value = "Cebollín"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 65
$ws.Range("K168").Value = 5000
$ws.Range("L168").Value = 5000
$ws.Range("M168").Value = 5000
$ws.Range("N168").Value = "`$/docena de paquetes"
$ws.Range("O168").Value = "Región de O'Higgins"
$ws.Range("P168").Value = 417
$ws.Range("Q168").Value = 12
$ws.Range("R168").Value = "Hortaliza"
